$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '22.434.55'
$ws.Range("E2").Value = '  -0.04%  '

# Row 3
$ws.Range("D3").Value = '1.568.07'
$ws.Range("E3").Value = '  +0.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("E5").Value = '  +0.06%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '288.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3702'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.83%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.30'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.37%  '

# Row 9
$ws.Range("E9").Value = '  -1.87%  '

# Row 10
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07490'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.08%  '

# Row 11
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.132'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.43%  '

# Row 12
$ws.Range("E12").Value = '  +0.09%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.57%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.923'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.97%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.866'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.15%  '

# Row 16
$ws.Range("D16").Value = '1.566.29'
$ws.Range("E16").Value = '  +0.31%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001116'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.52%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06752'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.35%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '87.50'
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = '  +0.02%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.337'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.30%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.26%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.21%  '

# Row 24
$ws.Range("D24").Value = '22.431.40'
$ws.Range("E24").Value = '  -0.03%  '

# Row 25
$ws.Range("E25").Value = '  -0.10%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.584'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.37%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.95%  '

# Row 28
$ws.Range("E28").Value = '  -0.45%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.013'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.19%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.30%  '

# Row 31
$ws.Range("D31").Value = '1.744.98'
$ws.Range("E31").Value = '  +0.48%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.053'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.39%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.011'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.23%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.110'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.51%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.752'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.60%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08371'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.01%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02464'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.20%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2257'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.08%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06403'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.284'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.72%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.331'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.23%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.30'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.13%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6302'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.45%  '

# Row 44
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.03%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.33%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6133'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.35%  '

# Row 47
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.775'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.27%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.057'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.35%  '

# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '126.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.51%  '

# Row 50
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.214'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.77%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07233'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.45%  '
